$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Text = "Speaker notes"
